$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 757, shifting the existing rows 757:798 down to 758:799
# (dimension grows from A1:D798 to A1:D799).
$ws.Rows.Item(757).Insert()

# Column A holds a date-like string (e.g. "2026/02/02") that must stay plain
# text, not get auto-converted into a date serial number. Force the cell to
# Text format before writing it, then clear the temporary formatting so the
# cell ends up with no explicit style, exactly like the other rows.
$ws.Range("A757").NumberFormat = "@"
$ws.Range("A757").Value = "2026/02/02"
$ws.Range("A757").ClearFormats()

$ws.Range("B757").Value = "月"
$ws.Range("C757").Value = 8
$ws.Range("D757").Value = 201
